$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (column-major order so new shared strings land
#     at the same indices the target workbook uses: 19..29) ---
$ws.Range("A2").Value = "controllermay02pipeline02@fpk12.com"
$ws.Range("A3").Value = "proctormay02pipeline02@fpk12.com"
$ws.Range("A4").Value = "examtakermay02pipeline02@fpk12.com"

$ws.Range("B2").Value = "masonsions0502"
$ws.Range("B3").Value = "jacksonions0502"
$ws.Range("B4").Value = "averyions0502"

$ws.Range("C2").Value = "harpersions0502"
$ws.Range("C3").Value = "jackisions0502"
$ws.Range("C4").Value = "wyattions0502"

$ws.Range("E3").Value = "LOC mar 0502"
$ws.Range("E4").Value = "LOC mar 0502"

# --- New row 9 ---
$ws.Range("E9").Value = "s"

# --- Column widths ---
# Column A: drop best-fit, set an explicit width (~35.22 chars). The
# engine quantizes ColumnWidth to 1/6-character steps, so feed it the
# value whose rounded result lands closest to the target width.
$ws.Columns("A").ColumnWidth = 34.385416666666664
# Column B -> 20 characters
$ws.Columns("B").ColumnWidth = 19.166666666666668
# Column C -> 16 characters
$ws.Columns("C").ColumnWidth = 15.166666666666666

# --- Selection ---
$ws.Range("N7").Select()
